$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on humidity(%) cells so Excel keeps them as text, not numeric percentages
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H46").NumberFormat = "@"

# Apply updated values
$ws.Range("E2").Value = "2026-02-13 05:18:27"
$ws.Range("N2").Value = "-3.0 °C 4:38 TU"
$ws.Range("O2").Value = "-2.1 °C"
$ws.Range("E3").Value = "2026-02-13 05:18:29"
$ws.Range("H3").Value = "67%"
$ws.Range("O3").Value = "-2.7 °C"
$ws.Range("E4").Value = "2026-02-13 05:18:32"
$ws.Range("H4").Value = "56%"
$ws.Range("J4").Value = "1003.1 hPa"
$ws.Range("N4").Value = "8.2 °C 4:48 TU"
$ws.Range("O4").Value = "9.9 °C"
$ws.Range("E5").Value = "2026-02-13 05:18:34"
$ws.Range("L5").Value = "29.2 km/h - 132º 4:50 TU"
$ws.Range("E6").Value = "2026-02-13 05:18:37"
$ws.Range("H6").Value = "61%"
$ws.Range("J6").Value = "1003.3 hPa"
$ws.Range("N6").Value = "7.5 °C 4:32 TU"
$ws.Range("E7").Value = "2026-02-13 05:18:39"
$ws.Range("J7").Value = "1003.9 hPa"
$ws.Range("N7").Value = "13.8 °C 4:34 TU"
$ws.Range("O7").Value = "14.6 °C"
$ws.Range("E8").Value = "2026-02-13 05:18:42"
$ws.Range("J8").Value = "1003.8 hPa"
$ws.Range("E9").Value = "2026-02-13 05:18:44"
$ws.Range("H9").Value = "64%"
$ws.Range("O9").Value = "8.6 °C"
$ws.Range("E10").Value = "2026-02-13 05:18:47"
$ws.Range("H10").Value = "75%"
$ws.Range("E11").Value = "2026-02-13 05:18:49"
$ws.Range("H11").Value = "78%"
$ws.Range("E12").Value = "2026-02-13 05:18:51"
$ws.Range("E13").Value = "2026-02-13 05:18:54"
$ws.Range("H13").Value = "83%"
$ws.Range("J13").Value = "1007.5 hPa"
$ws.Range("O13").Value = "-0.5 °C"
$ws.Range("E14").Value = "2026-02-13 05:18:56"
$ws.Range("H14").Value = "60%"
$ws.Range("N14").Value = "9.1 °C 4:59 TU"
$ws.Range("O14").Value = "11.7 °C"
$ws.Range("E15").Value = "2026-02-13 05:18:59"
$ws.Range("H15").Value = "65%"
$ws.Range("O15").Value = "9.0 °C"
$ws.Range("E16").Value = "2026-02-13 05:19:01"
$ws.Range("H16").Value = "59%"
$ws.Range("L16").Value = "69.8 km/h - 275º 4:57 TU"
$ws.Range("E17").Value = "2026-02-13 05:19:04"
$ws.Range("O17").Value = "1.5 °C"
$ws.Range("E18").Value = "2026-02-13 05:19:06"
$ws.Range("J18").Value = "1003.4 hPa"
$ws.Range("E19").Value = "2026-02-13 05:19:09"
$ws.Range("E20").Value = "2026-02-13 05:19:11"
$ws.Range("O20").Value = "-3.8 °C"
$ws.Range("E21").Value = "2026-02-13 05:19:13"
$ws.Range("H21").Value = "77%"
$ws.Range("J21").Value = "1006.1 hPa"
$ws.Range("O21").Value = "2.4 °C"
$ws.Range("E22").Value = "2026-02-13 05:19:16"
$ws.Range("H22").Value = "81%"
$ws.Range("E23").Value = "2026-02-13 05:19:18"
$ws.Range("H23").Value = "66%"
$ws.Range("L23").Value = "31.7 km/h - 142º 4:35 TU"
$ws.Range("M23").Value = "-2.1 °C 4:31 TU"
$ws.Range("O23").Value = "-3.6 °C"
$ws.Range("E24").Value = "2026-02-13 05:19:21"
$ws.Range("H24").Value = "84%"
$ws.Range("J24").Value = "1005.1 hPa"
$ws.Range("O24").Value = "6.7 °C"
$ws.Range("E25").Value = "2026-02-13 05:19:23"
$ws.Range("H25").Value = "57%"
$ws.Range("N25").Value = "-3.8 °C 4:30 TU"
$ws.Range("O25").Value = "-2.6 °C"
$ws.Range("E26").Value = "2026-02-13 05:19:26"
$ws.Range("J26").Value = "1003.8 hPa"
$ws.Range("O26").Value = "2.8 °C"
$ws.Range("E27").Value = "2026-02-13 05:19:29"
$ws.Range("H27").Value = "60%"
$ws.Range("M27").Value = "-0.7 °C 4:52 TU"
$ws.Range("E28").Value = "2026-02-13 05:19:31"
$ws.Range("H28").Value = "66%"
$ws.Range("J28").Value = "1003.8 hPa"
$ws.Range("E29").Value = "2026-02-13 05:19:33"
$ws.Range("H29").Value = "86%"
$ws.Range("K29").Value = "-0.1 MJ/m2"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-13 05:19:36"
$ws.Range("H30").Value = "74%"
$ws.Range("J30").Value = "1003.7 hPa"
$ws.Range("O30").Value = "7.6 °C"
$ws.Range("E31").Value = "2026-02-13 05:19:38"
$ws.Range("J31").Value = "1002.7 hPa"
$ws.Range("K31").Value = "-0.1 MJ/m2"
$ws.Range("E32").Value = "2026-02-13 05:19:41"
$ws.Range("H32").Value = "74%"
$ws.Range("N32").Value = "1.8 °C 4:47 TU"
$ws.Range("O32").Value = "5.0 °C"
$ws.Range("E33").Value = "2026-02-13 05:19:44"
$ws.Range("J33").Value = "1006.3 hPa"
$ws.Range("E34").Value = "2026-02-13 05:19:46"
$ws.Range("O34").Value = "-0.2 °C"
$ws.Range("E35").Value = "2026-02-13 05:19:49"
$ws.Range("J35").Value = "1005.9 hPa"
$ws.Range("O35").Value = "6.3 °C"
$ws.Range("E36").Value = "2026-02-13 05:19:51"
$ws.Range("H36").Value = "58%"
$ws.Range("J36").Value = "1003.3 hPa"
$ws.Range("N36").Value = "8.3 °C 4:43 TU"
$ws.Range("O36").Value = "11.8 °C"
$ws.Range("E37").Value = "2026-02-13 05:19:54"
$ws.Range("J37").Value = "1005.3 hPa"
$ws.Range("O37").Value = "4.1 °C"
$ws.Range("E38").Value = "2026-02-13 05:19:56"
$ws.Range("H38").Value = "51%"
$ws.Range("O38").Value = "10.8 °C"
$ws.Range("E39").Value = "2026-02-13 05:19:59"
$ws.Range("H39").Value = "49%"
$ws.Range("N39").Value = "-4.2 °C 4:50 TU"
$ws.Range("O39").Value = "-2.4 °C"
$ws.Range("E40").Value = "2026-02-13 05:20:01"
$ws.Range("J40").Value = "1007.3 hPa"
$ws.Range("E41").Value = "2026-02-13 05:20:04"
$ws.Range("H41").Value = "49%"
$ws.Range("J41").Value = "1004.4 hPa"
$ws.Range("K41").Value = "-0.1 MJ/m2"
$ws.Range("N41").Value = "10.1 °C 4:46 TU"
$ws.Range("O41").Value = "13.4 °C"
$ws.Range("E42").Value = "2026-02-13 05:20:06"
$ws.Range("O42").Value = "10.3 °C"
$ws.Range("E43").Value = "2026-02-13 05:20:09"
$ws.Range("H43").Value = "64%"
$ws.Range("N43").Value = "5.3 °C 4:46 TU"
$ws.Range("O43").Value = "7.7 °C"
$ws.Range("E44").Value = "2026-02-13 05:20:11"
$ws.Range("H44").Value = "79%"
$ws.Range("E45").Value = "2026-02-13 05:20:14"
$ws.Range("J45").Value = "1004.8 hPa"
$ws.Range("O45").Value = "2.4 °C"
$ws.Range("E46").Value = "2026-02-13 05:20:16"
$ws.Range("H46").Value = "80%"
$ws.Range("J46").Value = "1005.5 hPa"
$ws.Range("O46").Value = "7.1 °C"
